# Metrics.docx update
#   - Table 1 ("Metric / Tracking Method / Min / Max / Ave"):
#       "Number of Parameters" row, Max column: 5 -> 6
#   - Table 2 ("Metric / Week 3 / Week 4 / Week 5"), Week 5 column:
#       "# of Unit Tests"              -> 90
#       "Test Cases per Public Method" -> 90/56 = 1.61   (1.61 in red)
#       "Lines of Code per Method"     -> Max: 69    Ave: 6.8   (69 red, 6.8 green)
#       "Lines of Code"                -> 1878  (plus the _GoBack bookmark moves here)
#       "Number of Parameters"         -> Max: 6     Ave: 0.47 (both values green)

$d = $word.ActiveDocument

# Word RGB colors (r + g*256 + b*65536), matching the hex colors used in the diff.
$red   = 255        # FF0000
$green = 5287936     # 00B050

# ---------------------------------------------------------------------------
# 1) First table: "Number of Parameters" row, "Max" column: 5 -> 6
# ---------------------------------------------------------------------------
$t1 = $d.Tables(1)
$cell = $t1.Cell(6, 4)
$cell.Range.Find.Execute("5", $true, $false, $false, $false, $false, $true, 1, $false, "6", 2)

# ---------------------------------------------------------------------------
# 2) Second table, Week 5 column (column 4)
# ---------------------------------------------------------------------------
$t2 = $d.Tables(2)

# -- Row 2: "# of Unit Tests" -> 90 ------------------------------------------------
$cell = $t2.Cell(2, 4)
$cell.Range.Text = "90"

# -- Row 3: "Test Cases per Public Method" -> 90/56 = 1.61 (1.61 red) --------------
$cell = $t2.Cell(3, 4)
$r = $cell.Range
$cellStart = $r.Start
$r.Text = "90/56 = 1.61"

$redStart = $cellStart + 8
$redEnd = $redStart + 4
$redRange = $d.Range($redStart, $redEnd)
$redRange.Font.Color = $red

# -- Row 4: "Lines of Code per Method" -> Max: 69    Ave: 6.8 (69 red, 6.8 green) --
$cell = $t2.Cell(4, 4)
$r = $cell.Range
$cellStart = $r.Start
$r.Text = "Max: 69    Ave: 6.8"

$redStart = $cellStart + 5
$redEnd = $redStart + 2
$redRange = $d.Range($redStart, $redEnd)
$redRange.Font.Color = $red

$greenStart = $cellStart + 16
$greenEnd = $greenStart + 3
$greenRange = $d.Range($greenStart, $greenEnd)
$greenRange.Font.Color = $green

# -- Row 5: "Lines of Code" -> 1878 -------------------------------------------------
$cell = $t2.Cell(5, 4)
$cell.Range.Text = "1878"

# -- Row 6: "Number of Parameters" -> Max: 6     Ave: 0.47 (both values green) -----
$cell = $t2.Cell(6, 4)
$r = $cell.Range
$cellStart = $r.Start
$r.Text = "Max: 6     Ave: 0.47"

$greenStart = $cellStart + 5
$greenEnd = $greenStart + 2
$greenRange = $d.Range($greenStart, $greenEnd)
$greenRange.Font.Color = $green

$greenStart2 = $cellStart + 16
$greenEnd2 = $greenStart2 + 4
$greenRange2 = $d.Range($greenStart2, $greenEnd2)
$greenRange2.Font.Color = $green

# ---------------------------------------------------------------------------
# 3) Move the "_GoBack" bookmark from the "Number of Parameters"/Week 5 cell
#    to the end of the "Lines of Code"/Week 5 cell (best effort - this
#    bookmark simply marks the last edit location).
# ---------------------------------------------------------------------------
try {
    $old = $d.Bookmarks("_GoBack")
    $old.Delete()
} catch {
}

$cell = $t2.Cell(5, 4)
$r = $cell.Range
$r.Collapse(0)
try {
    $d.Bookmarks.Add("_GoBack", $r)
} catch {
}
